$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "Léptetőmotor + vezérlőegység" -> "Léptetőmotor + ULN2003A vezérlőegység"
#    split across three runs, with the _GoBack bookmark sitting between
#    the " ULN2003A" run and the " vezérlőegység" run (mirrors Word's
#    behaviour of re-anchoring _GoBack at the most recent edit point).
# ------------------------------------------------------------------

$target = $d.Content
$target.Find.ClearFormatting()
$null = $target.Find.Execute("Léptetőmotor + vezérlőegység")

# Shrink the found range down to just "Léptetőmotor +" so the trailing
# " vezérlőegység" text becomes a separate run when we insert into the gap.
$fullEnd = $target.End
$target.End = $target.Start + 15   # "Léptetőmotor +" = 15 characters

# Sanity anchor - remember the remainder of the original text.
$restStart = $target.End
$restEnd = $fullEnd

# Remove the old " vezérlőegység" tail and re-insert it after we've
# placed the new text, so every inserted chunk becomes its own run.
$restRange = $d.Range($restStart, $restEnd)
$restRange.Delete()

# Insert " ULN2003A" right after "Léptetőmotor +"
$insertPoint = $d.Range($target.End, $target.End)
$insertPoint.InsertAfter(" ULN2003A")
$insertPoint.Font.Name = "Courier New"
$insertPoint.Font.Size = 12

# Move the _GoBack bookmark here (delete old, add new) - this matches
# Word's automatic behaviour of tracking the most recent edit location.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$bmRange = $d.Range($insertPoint.End, $insertPoint.End)
$d.Bookmarks.Add("_GoBack", $bmRange)

# Insert " vezérlőegység" after the bookmark as its own run.
$tailPoint = $d.Range($insertPoint.End, $insertPoint.End)
$tailPoint.InsertAfter(" vezérlőegység")
$tailPoint.Font.Name = "Courier New"
$tailPoint.Font.Size = 12

# ------------------------------------------------------------------
# 2) Wrap "led" in a spell-check proofErr pair: " led" -> " " + "led"
#    (the proofErr elements themselves aren't part of the Word object
#    model text, but splitting "led" into its own run reproduces the
#    same run boundaries that Word's spell-checker produces.)
# ------------------------------------------------------------------

$ledRange = $d.Content
$null = $ledRange.Find.Execute("Error led")
$ledOnly = $d.Range($ledRange.End - 3, $ledRange.End)
$ledOnly.Select()
# Re-apply identical formatting to force a run split at this boundary.
$ledOnly.Font.Name = "Courier New"
$ledOnly.Font.Size = 12
